$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new "2022-Q3" summary row right after the header,
#    pushing the existing "2022-Q2" / "2021-Q2" rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# The row-insert copies formatting down from the row below (s="3" for
# B2:D2); the target has no explicit style there, so reset it.
$summary.Range("B2:D2").ClearFormats()

# A2 needs the same style as A3/A4 (s="2"); grab it from A3 (still has it).
$summary.Range("A3").Copy($summary.Range("A2"))

# New "2022-Q3" row.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.25

# Fix up the running index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" fund-holdings sheet. Duplicate the existing "2022-Q2"
#    sheet (same layout/header/styles), place it right before "2022-Q2",
#    rename it, then replace its data rows.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)   # "2022-Q2"
$q2.Copy($q2)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The source sheet has 7 data rows (2-8); the new sheet only needs 5 (2-6).
$q3.Rows.Item(8).Delete()
$q3.Rows.Item(7).Delete()

$fundData = @(
    @("003318", "景顺长城中证500行业中性低波动指数", "10.25", "93.67", "1.16", "0.1189", 4),
    @("001173", "中欧瑾和灵活配置混合 - A",           "2.26",  "92.00", "4.36", "0.0985", 9),
    @("512260", "华安中证500行业中性低波动ETF",        "1.07",  "97.91", "1.21", "0.0129", 4),
    @("001174", "中欧瑾和灵活配置混合 - C",           "0.23",  "92.00", "4.36", "0.0100", 9),
    @("501002", "长信价值优选混合",                    "0.39",  "93.83", "1.34", "0.0052", 10)
)

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $r = $i + 2
    $row = $fundData[$i]

    # Force the numeric-looking text columns (B..G) to stay text, then drop
    # the number-format style again so the cells end up styleless, matching
    # the source sheet's layout.
    $textRange = $q3.Range("B" + $r + ":G" + $r)
    $textRange.NumberFormat = "@"

    $q3.Range("A" + $r).Value = $i
    $q3.Range("B" + $r).Value = $row[0]
    $q3.Range("C" + $r).Value = $row[1]
    $q3.Range("D" + $r).Value = $row[2]
    $q3.Range("E" + $r).Value = $row[3]
    $q3.Range("F" + $r).Value = $row[4]
    $q3.Range("G" + $r).Value = $row[5]
    $q3.Range("H" + $r).Value = $row[6]

    $textRange.ClearFormats()
}
